$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price strings that look numeric (e.g. "303.41").
# Prefixing with an apostrophe forces Excel to keep them as literal text,
# matching the source data which stores these as text, not numbers.

$ws.Cells.Item(2, 4).Value = '43.935.69'
$ws.Cells.Item(2, 5).Value = '  +0.01%  '
$ws.Cells.Item(3, 4).Value = '2.231.31'
$ws.Cells.Item(3, 5).Value = '  -1.09%  '
$ws.Cells.Item(4, 5).Value = '  +0.18%  '
$ws.Cells.Item(5, 4).Value = '''303.41'
$ws.Cells.Item(5, 5).Value = '  -4.55%  '
$ws.Cells.Item(6, 4).Value = '''94.33'
$ws.Cells.Item(6, 5).Value = '  -7.06%  '
$ws.Cells.Item(7, 5).Value = '  -1.68%  '
$ws.Cells.Item(8, 5).Value = '  +0.23%  '
$ws.Cells.Item(9, 5).Value = '  -6.77%  '
$ws.Cells.Item(10, 4).Value = '''34.26'
$ws.Cells.Item(10, 5).Value = '  -7.74%  '
$ws.Cells.Item(11, 4).Value = '''0.0802'
$ws.Cells.Item(11, 5).Value = '  -3.53%  '
$ws.Cells.Item(12, 5).Value = '  -7.12%  '
$ws.Cells.Item(13, 5).Value = '  -2.96%  '
$ws.Cells.Item(14, 4).Value = '2.571.41'
$ws.Cells.Item(14, 5).Value = '  -1.26%  '
$ws.Cells.Item(15, 4).Value = '2.262.35'
$ws.Cells.Item(15, 5).Value = '  +0.09%  '
$ws.Cells.Item(16, 5).Value = '  -5.72%  '
$ws.Cells.Item(17, 4).Value = '''13.33'
$ws.Cells.Item(17, 5).Value = '  -7.90%  '
$ws.Cells.Item(18, 4).Value = '43.761.12'
$ws.Cells.Item(18, 5).Value = '  -0.18%  '
$ws.Cells.Item(19, 5).Value = '  -3.78%  '
$ws.Cells.Item(20, 4).Value = '''11.98'
$ws.Cells.Item(20, 5).Value = '  -11.35%  '
$ws.Cells.Item(21, 4).Value = '''6.11'
$ws.Cells.Item(21, 5).Value = '  -6.49%  '
$ws.Cells.Item(22, 4).Value = '''64.42'
$ws.Cells.Item(22, 5).Value = '  -2.08%  '
$ws.Cells.Item(23, 4).Value = '''235.73'
$ws.Cells.Item(23, 5).Value = '  +0.23%  '
$ws.Cells.Item(24, 5).Value = '  -7.35%  '
$ws.Cells.Item(26, 5).Value = '  -8.12%  '
$ws.Cells.Item(27, 4).Value = '''9.76'
$ws.Cells.Item(27, 5).Value = '  -3.90%  '
$ws.Cells.Item(28, 5).Value = '  -2.76%  '
$ws.Cells.Item(29, 4).Value = '''35.83'
$ws.Cells.Item(29, 5).Value = '  -3.34%  '
$ws.Cells.Item(30, 2).Value = 'EthereumClassic'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(30, 4).Value = '''19.86'
$ws.Cells.Item(30, 5).Value = '  -1.52%  '
$ws.Cells.Item(31, 2).Value = 'Filecoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(31, 4).Value = '''5.84'
$ws.Cells.Item(31, 5).Value = '  -6.15%  '
$ws.Cells.Item(32, 4).Value = '''152.46'
$ws.Cells.Item(32, 5).Value = '  -4.46%  '
$ws.Cells.Item(33, 4).Value = '''0.0801'
$ws.Cells.Item(33, 5).Value = '  -5.90%  '
$ws.Cells.Item(34, 4).Value = '''2.64'
$ws.Cells.Item(34, 5).Value = '  -2.15%  '
$ws.Cells.Item(35, 4).Value = '''3.23'
$ws.Cells.Item(35, 5).Value = '  +5.55%  '
$ws.Cells.Item(36, 2).Value = 'Stellar'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(36, 4).Value = '''0.117'
$ws.Cells.Item(36, 5).Value = '  -1.47%  '
$ws.Cells.Item(37, 2).Value = 'Kaspa'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(37, 4).Value = '''0.106'
$ws.Cells.Item(37, 5).Value = '  -7.57%  '
$ws.Cells.Item(38, 5).Value = '  -9.78%  '
$ws.Cells.Item(39, 4).Value = '''14.56'
$ws.Cells.Item(39, 5).Value = '  -9.92%  '
$ws.Cells.Item(40, 4).Value = '''3.78'
$ws.Cells.Item(40, 5).Value = '  -10.11%  '
$ws.Cells.Item(41, 5).Value = '  -11.66%  '
$ws.Cells.Item(42, 5).Value = '  -6.30%  '
$ws.Cells.Item(43, 5).Value = '  +0.23%  '
$ws.Cells.Item(44, 4).Value = '1.725.31'
$ws.Cells.Item(44, 5).Value = '  -5.14%  '
$ws.Cells.Item(45, 4).Value = '''83.63'
$ws.Cells.Item(45, 5).Value = '  +1.48%  '
$ws.Cells.Item(46, 5).Value = '  -6.97%  '
$ws.Cells.Item(47, 4).Value = '''98.89'
$ws.Cells.Item(47, 5).Value = '  -5.81%  '
$ws.Cells.Item(48, 4).Value = '''4.87'
$ws.Cells.Item(48, 5).Value = '  -6.74%  '
$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(49, 4).Value = '''14.37'
$ws.Cells.Item(49, 5).Value = '  +2.68%  '
$ws.Cells.Item(50, 2).Value = 'ordi'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Cells.Item(50, 4).Value = '''67.81'
$ws.Cells.Item(50, 5).Value = '  -10.51%  '
$ws.Cells.Item(51, 2).Value = 'FraxShare'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(51, 4).Value = '''7.93'
$ws.Cells.Item(51, 5).Value = '  -5.10%  '
